$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1,1,3,2,2,2,1,1,1,1,2,2,2,1,3,0,3,1,4,0,3,0,0,0,3,3,1,3,3,1,1,4,0,3,4,2,1,0,0,2,2,1,2,2,1,0,1,6,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}

Write-Output "Done updating G2:G50"
